# Fix: Allow re-signing a day (Upsert) instead of 400 Error.
# Ensures points recalculate if tasks are added late.
#
# The sheet used to carry 4 sample "planificacion" rows (rows 2-5). The
# fix re-seeds the sheet with a single, re-signed day (row 2 only):
#   - the date moves one day forward (re-signed),
#   - the ticket_id got corrupted/over-written with a huge numeric value,
#   - tecnico_nombre/patente/cliente were swapped for new placeholder
#     values ("Pedro Pascal" / "YRPT30" / "HOLA"),
#   - the old sample rows 3-5 (and the strings only they used) are gone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the old sample rows 3-5, keeping only the header + row 2 ---
$ws.Rows("3:5").Delete()

# --- Re-sign row 2: new date, corrupted ticket_id, new technician/plate/client ---
$ws.Range("A2").Value = 46022
$ws.Range("B2").Value = 33572342233334400
$ws.Range("I2").Value = "Pedro Pascal"
$ws.Range("J2").Value = "YRPT30"
$ws.Range("K2").Value = "HOLA"
$ws.Range("L2").Value = "GPS"
$ws.Range("M2").Value = " Corta Corriente"
$ws.Range("N2").Value = " Sensor Pta"
$ws.Range("O2").Value = " Sensor Temperatura"

# --- Column B now needs to fit the long literal ticket_id value ---
$ws.Columns.Item(2).ColumnWidth = 10.764322916666666

# --- Leave the selection where the author last left it ---
$ws.Range("K7").Select()
